# Roles.xlsx update — "Update biz plan, biz model"
#
# Builds out Sheet1 with the team-roles overview content: an intro
# paragraph, an "Overview" / "Geographic roles" header row (merged
# B2:E2, centered), a CEO row with the four geographic offices, and a
# "Functional roles include:" list. All populated cells use the larger
# (16pt) heading font used throughout this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1: intro / disclaimer text (col C) ---------------------------
$ws.Range("C1").Value = "This document indicates the roles of each member in the team. The role can be changed afterwards to fit with the condition of the company"
$ws.Range("C1").Font.Size = 16

# --- Row 2: section headers --------------------------------------------
$ws.Range("A2").Value = "Overview"
$ws.Range("A2").Font.Size = 16

$ws.Range("B2").Value = "Geographic roles"
$ws.Range("B2:E2").Font.Size = 16
$ws.Range("B2:E2").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B2:E2").Merge()

# --- Column A: CEO, then the functional-roles list (rows 3-13) --------
$ws.Range("A3").Value = "CEO"
$ws.Range("A3").Font.Size = 16

$ws.Range("A4").Value = "Functional roles include:"
$ws.Range("A4").Font.Size = 16

$ws.Range("A5").Value = "Area manager"
$ws.Range("A5").Font.Size = 16

$ws.Range("A6").Value = "Finance"
$ws.Range("A6").Font.Size = 16

$ws.Range("A7").Value = "Marketing"
$ws.Range("A7").Font.Size = 16

$ws.Range("A8").Value = "Sales"
$ws.Range("A8").Font.Size = 16

$ws.Range("A9").Value = "Negotiation"
$ws.Range("A9").Font.Size = 16

$ws.Range("A10").Value = "Production"
$ws.Range("A10").Font.Size = 16

$ws.Range("A11").Value = "Logistics"
$ws.Range("A11").Font.Size = 16

$ws.Range("A12").Value = "R&D"
$ws.Range("A12").Font.Size = 16

$ws.Range("A13").Value = "Etc."
$ws.Range("A13").Font.Size = 16

# --- Row 3, columns B-E: the geographic offices for the CEO row -------
$ws.Range("B3").Value = "Home Office (Leich.)"
$ws.Range("B3").Font.Size = 16

$ws.Range("C3").Value = "US"
$ws.Range("C3").Font.Size = 16

$ws.Range("D3").Value = "EU"
$ws.Range("D3").Font.Size = 16

$ws.Range("E3").Value = "Brazil"
$ws.Range("E3").Font.Size = 16

# --- Column widths for the label / office columns ----------------------
$ws.Columns.Item(1).ColumnWidth = 27.25
$ws.Columns.Item(2).ColumnWidth = 22.75

# --- Page setup ----------------------------------------------------------
$ws.PageSetup.Orientation = 1   # xlPortrait

# --- Selection, matching the saved workbook's cursor position ---------
$ws.Range("A6").Select()
